{"js": "const replacements = [\n  [\"2024-09-03 Tuesday\", \"2024-09-04 Wednesday\"],\n  [\"113\u00f77=\", \"412\u00f78=\"],\n  [\"500\u00f74=\", \"120\u00f76=\"],\n  [\"678\u00f73=\", \"145\u00f72=\"],\n  [\"145\u00f76=\", \"262\u00f72=\"],\n  [\"400\u00f72=\", \"646\u00f74=\"],\n  [\"550\u00f74=\", \"615\u00f74=\"],\n  [\"510\u00f75=\", \"342\u00f74=\"],\n  [\"816\u00f72=\", \"906\u00f76=\"],\n  [\"664\u00f74=\", \"791\u00f78=\"],\n  [\"334\u00f79=\", \"310\u00f77=\"],\n  [\"706\u00f75=\", \"883\u00f76=\"],\n  [\"383\u00f79=\", \"449\u00f78=\"],\n  [\"856\u00f74=\", \"366\u00f75=\"],\n  [\"176\u00f76=\", \"409\u00f79=\"],\n  [\"778\u00f79=\", \"273\u00f78=\"],\n  [\"370\u00f79=\", \"771\u00f72=\"],\n  [\"119\u00f78=\", \"242\u00f78=\"],\n  [\"355\u00f78=\", \"665\u00f76=\"],\n  [\"158\u00f78=\", \"139\u00f79=\"],\n  [\"962\u00f72=\", \"752\u00f78=\"],\n  [\"761\u00f77=\", \"407\u00f78=\"],\n  [\"519\u00f73=\", \"857\u00f74=\"],\n  [\"595\u00f73=\", \"661\u00f76=\"],\n  [\"350\u00f72=\", \"726\u00f72=\"],\n  [\"490\u00f78=\", \"260\u00f78=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-09-03 Tuesday\", \"2024-09-04 Wednesday\"),\n  @(\"113\u00f77=\", \"412\u00f78=\"),\n  @(\"500\u00f74=\", \"120\u00f76=\"),\n  @(\"678\u00f73=\", \"145\u00f72=\"),\n  @(\"145\u00f76=\", \"262\u00f72=\"),\n  @(\"400\u00f72=\", \"646\u00f74=\"),\n  @(\"550\u00f74=\", \"615\u00f74=\"),\n  @(\"510\u00f75=\", \"342\u00f74=\"),\n  @(\"816\u00f72=\", \"906\u00f76=\"),\n  @(\"664\u00f74=\", \"791\u00f78=\"),\n  @(\"334\u00f79=\", \"310\u00f77=\"),\n  @(\"706\u00f75=\", \"883\u00f76=\"),\n  @(\"383\u00f79=\", \"449\u00f78=\"),\n  @(\"856\u00f74=\", \"366\u00f75=\"),\n  @(\"176\u00f76=\", \"409\u00f79=\"),\n  @(\"778\u00f79=\", \"273\u00f78=\"),\n  @(\"370\u00f79=\", \"771\u00f72=\"),\n  @(\"119\u00f78=\", \"242\u00f78=\"),\n  @(\"355\u00f78=\", \"665\u00f76=\"),\n  @(\"158\u00f78=\", \"139\u00f79=\"),\n  @(\"962\u00f72=\", \"752\u00f78=\"),\n  @(\"761\u00f77=\", \"407\u00f78=\"),\n  @(\"519\u00f73=\", \"857\u00f74=\"),\n  @(\"595\u00f73=\", \"661\u00f76=\"),\n  @(\"350\u00f72=\", \"726\u00f72=\"),\n  @(\"490\u00f78=\", \"260\u00f78=\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
